$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells F1:K1
$ws.Range("F1").Value = "C4.5 acc"
$ws.Range("G1").Value = "credal-C4.5 acc"
$ws.Range("H1").Value = "SPN acc"
$ws.Range("I1").Value = "CSPN low"
$ws.Range("J1").Value = "CSPN high"
$ws.Range("K1").Value = "CSPN robust"

# Copy the style of an existing header cell (A1) into F1:K1
$ws.Range("A1").Copy()
$ws.Range("F1:K1").PasteSpecial(-4122) # xlPasteFormats

# Updated values for B2:E6
$ws.Range("B2").Value = 86.93957115009746
$ws.Range("C2").Value = 78.55750487329435
$ws.Range("D2").Value = 93.78167641325535
$ws.Range("E2").Value = 92.64827164734079

$ws.Range("B3").Value = 87.34892787524366
$ws.Range("C3").Value = 78.7719298245614
$ws.Range("D3").Value = 93.99610136452242
$ws.Range("E3").Value = 92.90819123574961

$ws.Range("B4").Value = 87.83625730994153
$ws.Range("C4").Value = 78.61598440545808
$ws.Range("D4").Value = 94.21052631578947
$ws.Range("E4").Value = 93.13476415563333

$ws.Range("B5").Value = 86.93957115009746
$ws.Range("C5").Value = 76.80311890838206
$ws.Range("D5").Value = 93.56725146198829
$ws.Range("E5").Value = 92.25413876075127

$ws.Range("B6").Value = 85.61403508771929
$ws.Range("C6").Value = 74.56140350877193
$ws.Range("D6").Value = 94.03508771929823
$ws.Range("E6").Value = 92.5942539809449

# New data for F2:K6
$ws.Range("F2").Value = 86.82261208576999
$ws.Range("G2").Value = 87.54385964912281
$ws.Range("H2").Value = 90.50682261208576
$ws.Range("I2").Value = 90.50682261208576
$ws.Range("J2").Value = 90.50682261208576
$ws.Range("K2").Value = 90.50682261208576

$ws.Range("F3").Value = 81.2280701754386
$ws.Range("G3").Value = 86.78362573099415
$ws.Range("H3").Value = 90.11695906432747
$ws.Range("I3").Value = 90.35087719298245
$ws.Range("J3").Value = 90.38986354775828
$ws.Range("K3").Value = 90.38585024653136

$ws.Range("F4").Value = 62.90448343079922
$ws.Range("G4").Value = 86.70565302144249
$ws.Range("H4").Value = 90.29239766081872
$ws.Range("I4").Value = 90.9551656920078
$ws.Range("J4").Value = 91.05263157894736
$ws.Range("K4").Value = 91.04391698199747

$ws.Range("F5").Value = 62.35867446393762
$ws.Range("G5").Value = 84.9317738791423
$ws.Range("H5").Value = 90.33138401559454
$ws.Range("I5").Value = 90.5653021442495
$ws.Range("J5").Value = 90.60428849902534
$ws.Range("K5").Value = 90.60130718954248

$ws.Range("F6").Value = 56.35477582846003
$ws.Range("G6").Value = 82.70955165692007
$ws.Range("H6").Value = 89.39571150097466
$ws.Range("I6").Value = 89.20077972709552
$ws.Range("J6").Value = 89.23976608187134
$ws.Range("K6").Value = 89.23529411764706
